$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.733529806137085
$ws.Range("B1").Value = 2.841336011886597
$ws.Range("C1").Value = 3.50223708152771
$ws.Range("D1").Value = 1.352825164794922
$ws.Range("E1").Value = 0.9037905931472778
